# Fruta / hortaliza, semanal
# Insert a new daily price record as row 42 on the "Naranja" sheet,
# pushing the existing rows 42-135 down to 43-136 (dimension grows from
# A1:T135 to A1:T136).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 42 - this shifts every row at/after
# 42 down by one, automatically extending the sheet dimension.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new record's data.
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C42").Value = "Arica y Parinacota"
$ws.Range("D42").Value = 45014
$ws.Range("E42").Value = 15
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100102
$ws.Range("H42").Value = "Cítricos"
$ws.Range("I42").Value = 100102005
$ws.Range("J42").Value = "Naranja"
$ws.Range("K42").Value = "Valencia"
$ws.Range("L42").Value = "Tercera"
$ws.Range("M42").Value = 230
$ws.Range("N42").Value = 1100
$ws.Range("O42").Value = 1150
$ws.Range("P42").Value = 1117
$ws.Range("Q42").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R42").Value = "Región de Coquimbo"
$ws.Range("S42").Value = 1117
$ws.Range("T42").Value = 1
